# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the Q1-2022 fund-holdings breakdown (same layout as
#    the other quarterly sheets).
# 2. Insert a new row at the top of the "总计" (summary) sheet's data for
#    the 2022-Q1 aggregate (date / count / value), pushing the older rows
#    down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: insert the "2022-Q1" worksheet before "总计"
# ---------------------------------------------------------------------
$sheets = $wb.Worksheets
$totalSheet = $sheets.Item("总计")
$newSheet = $sheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Reference sheet that already has the identical column layout/styling
# ("基金代码","基金名称","基金规模", ...) so we can clone header/column-A
# formatting (bold + border = style used on row 1 and column A).
$refSheet = $sheets.Item("2021-Q4")

# --- header row ---------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2  # column B
foreach ($h in $headers) {
    $refSheet.Range("B1").Copy()
    $newSheet.Cells.Item(1, $col).PasteSpecial(-4122)
    $newSheet.Cells.Item(1, $col).Value = $h
    $col++
}

# --- data rows ------------------------------------------------------
# Each row: fund code, fund name, fund size, stock position, position
# ratio, held market value (100M CNY), position rank.
# All but the first (id) and last (rank) columns are stored as TEXT in
# the source workbook (e.g. "49.29"), matching the other quarter sheets.
$rows = @(
    @("002910", "易方达供给改革灵活配置混合", "49.29", "87.54", "6.20", "3.0560", 4),
    @("070021", "嘉实主题新动力混合", "24.04", "93.93", "4.54", "1.0914", 10),
    @("000985", "嘉实逆向策略股票", "13.64", "93.90", "4.55", "0.6206", 10),
    @("206002", "鹏华精选成长混合", "4.48", "92.68", "4.23", "0.1895", 9),
    @("008132", "鹏华价值驱动混合", "4.33", "90.91", "3.30", "0.1429", 8),
    @("004818", "国寿安保目标策略灵活配置混合A", "4.06", "36.45", "2.20", "0.0893", 6),
    @("014307", "嘉实多元动力混合A", "1.83", "91.81", "4.60", "0.0842", 6),
    @("004819", "国寿安保目标策略灵活配置混合C", "0.57", "36.45", "2.20", "0.0125", 6),
    @("014308", "嘉实多元动力混合C", "0.20", "91.81", "4.60", "0.0092", 6)
)

$r = 2
foreach ($row in $rows) {
    $fundCode = $row[0]
    $fundName = $row[1]
    $fundSize = $row[2]
    $stockPos = $row[3]
    $posRatio = $row[4]
    $heldValue = $row[5]
    $posRank = $row[6]

    # Column A: running index (0-based), bold+bordered style like the
    # other quarter sheets.
    $refSheet.Range("A2").Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $newSheet.Cells.Item($r, 1).Value = ($r - 2)

    # Column B: fund code -- force text so leading zeros survive.
    $newSheet.Cells.Item($r, 2).Value = "'" + $fundCode
    $newSheet.Cells.Item($r, 2).ClearFormats()

    # Column C: fund name (plain text, no numeric look-alike issue).
    $newSheet.Cells.Item($r, 3).Value = $fundName

    # Column D: fund size -- numeric-looking text.
    $newSheet.Cells.Item($r, 4).Value = "'" + $fundSize
    $newSheet.Cells.Item($r, 4).ClearFormats()

    # Column E: stock position -- numeric-looking text.
    $newSheet.Cells.Item($r, 5).Value = "'" + $stockPos
    $newSheet.Cells.Item($r, 5).ClearFormats()

    # Column F: position ratio -- numeric-looking text.
    $newSheet.Cells.Item($r, 6).Value = "'" + $posRatio
    $newSheet.Cells.Item($r, 6).ClearFormats()

    # Column G: held market value -- numeric-looking text.
    $newSheet.Cells.Item($r, 7).Value = "'" + $heldValue
    $newSheet.Cells.Item($r, 7).ClearFormats()

    # Column H: position rank -- a real number.
    $newSheet.Cells.Item($r, 8).Value = $posRank

    $r++
}

# ---------------------------------------------------------------------
# Part 2: insert the 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Column A keeps the running index + bold/bordered style (copy from the
# row that used to be row 2, now shifted to row 3).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("B2").ClearFormats()

$totalSheet.Range("C2").Value = 9
$totalSheet.Range("C2").ClearFormats()

$totalSheet.Range("D2").Value = 5.3
$totalSheet.Range("D2").ClearFormats()

# The shifted-down rows (old A2:A6 => new A3:A7) kept their original
# 0-based index values; bump each by one so the running index stays
# contiguous (0,1,2,3,4,5) after the new row was spliced in at the top.
for ($row = 7; $row -ge 3; $row--) {
    $oldIdx = $totalSheet.Cells.Item($row, 1).Value2
    $totalSheet.Cells.Item($row, 1).Value2 = $oldIdx + 1
}
